$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 (Caso 6417, NUMANCIA 436) - entire row shifts everything up
$ws.Rows.Item(6).Delete()

# After the above delete, old row 9 (ARTIGAS, JOSE GERVASIO 924) is now row 8. Delete it too.
$ws.Rows.Item(8).Delete()

# Remaining data rows are now, in order:
#  6 = old row7  (ALBERDI, JUAN BAUTISTA AV. 1091 / Caso 6557)
#  7 = old row8  (POLA 591 / Caso 6193)
#  8 = old row10 (Carlos E. Pellegrini 6030 / Caso 6568)

# --- Update row 6 (ALBERDI): only the Observaciones (H) text changes ---
$ws.Range("H6").Value = "Verificar si la linga panseada es nuestra y ver con pablo como resolverlo"

# --- Row 7 becomes the Pellegrini case (Caso 6568), with OT updated to ICD30313541 ---
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "6568"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "8/5/2025"
$ws.Range("C7").Value = "Carlos E. Pellegrini 6030"
$ws.Range("D7").Value = 12
$ws.Range("E7").Value = "ICD30313541 "
$ws.Range("H7").Value = "Tendido a muy baja altura se solicita retiro o levantarlo"
$ws.Range("J7").Value = '{"direccionesNormalizadas": [{"altura": 6030, "cod_calle": 17053, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.504789", "y": "-34.564505"}, "direccion": "PELLEGRINI, CARLOS E. 6030, CABA", "nombre_calle": "PELLEGRINI, CARLOS E.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K7").Value = -58.504789
$ws.Range("L7").Value = -34.564505
$ws.Range("M7").Value = "Paternal"
$ws.Range("N7").Value = "Capital Norte"

# --- Row 8 becomes the new CAMPANA 382 case (Caso 6569) ---
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "6569"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "8/6/2025"
$ws.Range("C8").Value = "CAMPANA 382"
$ws.Range("D8").Value = 10
$ws.Range("E8").Value = "Pendiente de ADM"
$ws.Range("H8").Value = "Tendido a baja altura"
$ws.Range("J8").Value = '{"direccionesNormalizadas": [{"altura": 382, "cod_calle": 3039, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.476505", "y": "-34.628022"}, "direccion": "CAMPANA 382, CABA", "nombre_calle": "CAMPANA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K8").Value = -58.476505
$ws.Range("L8").Value = -34.628022
$ws.Range("M8").Value = "Devoto"
$ws.Range("N8").Value = "Capital Norte"
